$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 holds the "(Leon Thomm)" team member entry that is being removed
# from the roster. Deleting the whole row shifts everything below it up by
# one, which also re-targets the SUM(...)/B5*100-B.. formulas so they keep
# covering the (now shorter) list of team members.
$ws.Rows.Item(13).Delete()

# The conditional formatting that used to flag a blank C13 ("hours" for the
# removed member) no longer applies to anything, so drop that rule.
$fcC13 = $ws.Range("C13").FormatConditions.Item(1)
$fcC13.Delete()

# The two "fill in your tasks" column rules used to cover A20:A26/B20:B26;
# after the row shift they must cover A19:A25/B19:B25 instead.
$fcA = $ws.Range("A20:A26").FormatConditions.Item(1)
$fcA.ModifyAppliesToRange($ws.Range("A19:A25"))
$fcA.Formula1 = "=ISBLANK(`$A`$19)"

$fcB = $ws.Range("B20:B26").FormatConditions.Item(1)
$fcB.ModifyAppliesToRange($ws.Range("B19:B25"))
$fcB.Formula1 = "=ISBLANK(`$B`$19)"

# Update the active selection to match the post-edit state.
$ws.Range("D19").Select()
